# Updated cryptos list with GitHub Actions
# Applies refreshed Price (D) and Volume(1h) (E) values for rows 2-51.
# Cells whose new Price text parses as a plain number get NumberFormat "@"
# set first so Excel keeps them as text (matching the original inlineStr
# string cells) instead of silently converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.421.93"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").Value = "3.766.50"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.63"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.77"
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("D7").Value = "3.765.55"
$ws.Range("E7").Value = "  -0.99%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("E11").Value = "  +4.74%  "

$ws.Range("E12").Value = "  -1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.18"
$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  -3.62%  "

$ws.Range("D15").Value = "4.391.42"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").Value = "3.766.62"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").Value = "69.502.36"
$ws.Range("E17").Value = "  -0.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.55"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("E19").Value = "  -3.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.21"
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.36"
$ws.Range("E21").Value = "  -2.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.39"
$ws.Range("E22").Value = "  -1.94%  "

$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.50"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.88"
$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("E27").Value = "  -3.89%  "

$ws.Range("E28").Value = "  -5.43%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  +2.70%  "

$ws.Range("E33").Value = "  -2.82%  "

$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("E38").Value = "  +3.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.341"
$ws.Range("E39").Value = "  +2.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "456.81"
$ws.Range("E40").Value = "  +9.20%  "

$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.81"
$ws.Range("E42").Value = "  -2.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("E43").Value = "  +5.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.83"
$ws.Range("E44").Value = "  -1.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.59"
$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("D46").Value = "2.958.97"
$ws.Range("E46").Value = "  -2.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0361"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.42"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.13"

$ws.Range("E51").Value = "  +1.02%  "
